# Updates the crypto price/volume table (and two swapped rows) per the
# latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.830.23'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '2.344.20'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '239.45'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '0.667'
$ws.Range('E6').Value = '  -2.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '72.39'
$ws.Range('E7').Value = '  -4.03%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.589'
$ws.Range('E9').Value = '  -5.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '0.0995'
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '58.43'
$ws.Range('E11').Value = '  +2.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '32.47'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '7.15'
$ws.Range('E14').Value = '  -4.02%  '
$ws.Range('D15').Value = '2.693.76'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '16.23'
$ws.Range('E16').Value = '  -3.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '0.896'
$ws.Range('E17').Value = '  -2.05%  '
$ws.Range('D18').Value = '2.342.42'
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('D19').Value = '43.735.30'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '6.65'
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '77.80'
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '253.19'
$ws.Range('E23').Value = '  -0.74%  '
$ws.Range('E24').Value = '  +10.07%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '3.73'
$ws.Range('E26').Value = '  +2.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '2.49'
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '10.36'
$ws.Range('E28').Value = '  -5.16%  '
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '175.91'
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '22.24'
$ws.Range('E31').Value = '  -3.19%  '
$ws.Range('E32').Value = '  -0.85%  '
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '0.0738'
$ws.Range('E34').Value = '  -1.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '5.08'
$ws.Range('E35').Value = '  -3.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '5.36'
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '2.38'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('B39').Value = 'THORChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '6.37'
$ws.Range('E39').Value = '  -1.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '0.0272'
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '5.24'
$ws.Range('E41').Value = '  +16.69%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '63.90'
$ws.Range('E42').Value = '  +17.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '9.15'
$ws.Range('E43').Value = '  +2.73%  '
$ws.Range('E44').Value = '  +6.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '18.70'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '0.197'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '2.45'
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '1.22'
$ws.Range('E49').Value = '  -2.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '98.20'
$ws.Range('E50').Value = '  -3.25%  '
$ws.Range('E51').Value = '  -4.22%  '
